# Update evaluations on the QuantitativeMetrics sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# Runtime without error: yes -> no
$ws.Range("B6").Value = "no"

# Assertion validity row: clear the yes/note values (now blank / numeric-empty)
$ws.Range("B7").Value = $null
$ws.Range("C7").Value = $null

# Code BLEU score update
$ws.Range("B12").Value = 0.3253501897437739
$ws.Range("C12").Value = "{'codebleu': 0.32535018974377394, 'ngram_match_score': 0.06773015082000793, 'weighted_ngram_match_score': 0.1054654799499596, 'syntax_match_score': 0.6043956043956044, 'dataflow_match_score': 0.5238095238095238}"

# Update the active selection on the sheet to B7 (matches saved cursor position)
$ws.Range("B7").Select()
